$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.708.11'
$ws.Range('D3').Value = '2.060.27'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  -0.07%  '
$c = $ws.Range('D5')
$c.Formula = "'244.48"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '
$c = $ws.Range('D6')
$c.Formula = "'0.668"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('E7').Value = '  +0.02%  '
$c = $ws.Range('D8')
$c.Formula = "'55.41"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -6.36%  '
$c = $ws.Range('D9')
$c.Formula = "'60.55"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.60%  '
$c = $ws.Range('D11')
$c.Formula = "'0.0752"
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -3.05%  '
$ws.Range('E12').Value = '  -3.03%  '
$ws.Range('E13').Value = '  +7.41%  '
$c = $ws.Range('D14')
$c.Formula = "'14.85"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -4.30%  '
$ws.Range('D15').Value = '2.361.19'
$ws.Range('E15').Value = '  +0.08%  '
$c = $ws.Range('D16')
$c.Formula = "'5.51"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -3.88%  '
$ws.Range('D17').Value = '2.068.17'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '36.639.60'
$ws.Range('E18').Value = '  -1.03%  '
$c = $ws.Range('D19')
$c.Formula = "'17.42"
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -4.20%  '
$c = $ws.Range('D20')
$c.Formula = "'72.37"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -2.44%  '
$ws.Range('D21').Value = '0.0₃0866'
$ws.Range('E21').Value = '  -2.91%  '
$c = $ws.Range('D22')
$c.Formula = "'238.86"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  -3.77%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  -2.62%  '
$c = $ws.Range('D26')
$c.Formula = "'2.26"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +5.03%  '
$c = $ws.Range('D27')
$c.Formula = "'9.28"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -7.16%  '
$c = $ws.Range('D28')
$c.Formula = "'166.11"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -1.90%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('E31').Value = '  +7.56%  '
$c = $ws.Range('D32')
$c.Formula = "'5.11"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -6.33%  '
$c = $ws.Range('D33')
$c.Formula = "'4.53"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -3.82%  '
$c = $ws.Range('D34')
$c.Formula = "'0.0599"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -3.04%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('E37').Value = '  +1.19%  '
$c = $ws.Range('D38')
$c.Formula = "'2.23"
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -4.21%  '
$c = $ws.Range('D39')
$c.Formula = "'5.11"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -3.29%  '
$c = $ws.Range('D40')
$c.Formula = "'1.26"
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -5.83%  '
$ws.Range('E41').Value = '  -5.94%  '
$ws.Range('E42').Value = '  -3.36%  '
$ws.Range('E43').Value = '  -5.01%  '
$c = $ws.Range('D44')
$c.Formula = "'95.19"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -2.92%  '
$c = $ws.Range('D45')
$c.Formula = "'0.0919"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -4.78%  '
$ws.Range('D46').Value = '1.421.43'
$ws.Range('E46').Value = '  +9.04%  '
$c = $ws.Range('D47')
$c.Formula = "'7.64"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +11.90%  '
$c = $ws.Range('D48')
$c.Formula = "'16.05"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -6.05%  '
$c = $ws.Range('D49')
$c.Formula = "'2.93"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.83%  '
$c = $ws.Range('D50')
$c.Formula = "'2.28"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -3.60%  '
$ws.Range('D51').Value = '2.246.26'
$ws.Range('E51').Value = '  -0.03%  '
